$d = $word.ActiveDocument

function Set-ParagraphXml($para, [string]$xml) {
    $r = $para.Range
    $r.InsertXML($xml)
}

# --- 1) "Dayofweek()" heading paragraph: split "Dayofweek()" into two runs with
#        spell-check markers, and split the "select dayofweek(...)" example text
#        so "dayofweek" is separately spell-marked. ---
$xmlDayofweek = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="234AEA61" w14:textId="28A2ECB2" w:rsidR="00FB762C" w:rsidRDefault="00FB762C"><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="006D7AD2"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Dayofweek</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="006D7AD2"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>()</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="004D7507"><w:t>–</w:t></w:r><w:r w:rsidR="006D7AD2"><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Ex</w:t></w:r><w:r w:rsidR="004D7507"><w:t xml:space="preserve">: </w:t></w:r><w:r w:rsidR="004D7507" w:rsidRPr="004D7507"><w:t xml:space="preserve">select </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="004D7507" w:rsidRPr="004D7507"><w:t>dayofweek</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="004D7507" w:rsidRPr="004D7507"><w:t>("1993-02-04") as "my birth day";</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$found = $d.Content
$found.Find.ClearFormatting()
$found.Find.Execute("Dayofweek() ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p1 = $found.Paragraphs(1)
Set-ParagraphXml $p1 $xmlDayofweek

# --- 2) "Concat()" heading paragraph (first SELECT CONCAT example). ---
$xmlConcatHeading = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="0BC3F57B" w14:textId="4CF6DA3E" w:rsidR="00FB762C" w:rsidRDefault="00FB762C" w:rsidP="004D7507"><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="006D7AD2"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Concat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="006D7AD2"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>()</w:t></w:r><w:r w:rsidR="006D7AD2"><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">- Ex- </w:t></w:r><w:r w:rsidR="004D7507"><w:t>SELECT CONCAT(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="004D7507"><w:t>first_name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="004D7507"><w:t xml:space="preserve">, " " , </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="004D7507"><w:t>last_name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="004D7507"><w:t xml:space="preserve">) AS </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="004D7507"><w:t>full_name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="004D7507"><w:t xml:space="preserve"> FROM actor;</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$found2 = $d.Content
$found2.Find.ClearFormatting()
$found2.Find.Execute("Concat() ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p2 = $found2.Paragraphs(1)
Set-ParagraphXml $p2 $xmlConcatHeading

# --- 3) Body paragraph containing the second "SELECT CONCAT(...)" occurrence. ---
$xmlConcatBody = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="1FDA9293" w14:textId="7105F387" w:rsidR="00C621AE" w:rsidRDefault="00FB762C" w:rsidP="006D7AD2"><w:r><w:t xml:space="preserve">This would be very useful for displaying two pieces of information together. Like a name, address, or name and phone number. </w:t></w:r><w:r w:rsidR="006D7AD2"><w:t xml:space="preserve">Please see pictures of this execution numbered on the following page. </w:t></w:r><w:r w:rsidR="007465FE"><w:t>For this example, I had to do some digging to find some information I could concatenate from a table. I didn’t want to bother with only strings. First</w:t></w:r><w:r w:rsidR="006D7AD2"><w:t>,</w:t></w:r><w:r w:rsidR="007465FE"><w:t xml:space="preserve"> I ran SHOW DATABASES; Then I selected USE SAKILA;</w:t></w:r><w:r w:rsidR="001B732F"><w:t xml:space="preserve"> (</w:t></w:r><w:r w:rsidR="006D7AD2"><w:t>#</w:t></w:r><w:r w:rsidR="001B732F"><w:t xml:space="preserve"> 1) as I saw this imported as an example database in my download. At first, I tried to use SHOW COLUMNS; but I realized that I wasn’t in a table yet. I ran SHOW TABLES; then I selected actor from the list of tables and ran SHOW COLUMNS FROM actor; (</w:t></w:r><w:r w:rsidR="006D7AD2"><w:t>#</w:t></w:r><w:r w:rsidR="001B732F"><w:t xml:space="preserve"> 2). Now I see that I have some names I can concatenate.</w:t></w:r><w:r w:rsidR="006D7AD2"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="006D7AD2"><w:t>SELECT CONCAT(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006D7AD2"><w:t>first_name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006D7AD2"><w:t xml:space="preserve">, " " , </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006D7AD2"><w:t>last_name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006D7AD2"><w:t xml:space="preserve">) AS </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006D7AD2"><w:t>full_name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="006D7AD2"><w:t xml:space="preserve"> FROM actor;</w:t></w:r><w:r w:rsidR="006D7AD2"><w:t xml:space="preserve"> this query returned 200 rows of names (# 3 &amp; 4). </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$found3 = $d.Content
$found3.Find.ClearFormatting()
$found3.Find.Execute("this query returned 200 rows", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p3 = $found3.Paragraphs(1)
Set-ParagraphXml $p3 $xmlConcatBody

Write-Output "Done"
